$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in missing B values for existing rows 10-12 ---
$ws.Range("B10").Value = 0.54894211041483298
$ws.Range("B11").Value = 0.55229771608690903
$ws.Range("B12").Value = 0.55010864533221904

# --- Row 13: swap B13 and C13 values ---
$ws.Range("B13").Value = 0.54667524152360103
$ws.Range("C13").Value = 0.40581239595836999

# --- New rows 14-23 ---
# NOTE: shared-string pool is built in first-seen order, and the target
# workbook has string index 13 = "...Reg=0.2256" (row 15) and index 14 =
# "...Reg=0.2256, RandomSeed=421" (row 14). So we register the row-15
# string first, then the row-14 string, before assigning them to cells in
# row order.
$ws.Range("A15").Value = "Elopp V1, 200 iterations, Reg=0.2256"
$ws.Range("A14").Value = "Elopp V1, 200 iterations, Reg=0.2256, RandomSeed=421"

$ws.Range("B14").Value = 0.548950207783092
$ws.Range("C14").Value = 0.406200483828216

$ws.Range("B15").Value = 0.54946373351655198
$ws.Range("C15").Value = 0.40645411249053998

$ws.Range("A16").Value = "Elopp V1, 400 iterations, Reg=0.2256"
$ws.Range("B16").Value = 0.54991849851126096
$ws.Range("C16").Value = 0.40667455430965399

$ws.Range("A17").Value = "Elopp V1, 400 iterations, Reg=0.2256, RandomSeed=421"
$ws.Range("B17").Value = 0.54979899562884904
$ws.Range("C17").Value = 0.406682990358605

$ws.Range("A18").Value = "Elopp V1, 600 iterations, Reg=0.2256, RandomSeed=421"
$ws.Range("B18").Value = 0.54974746327577895
$ws.Range("C18").Value = 0.40660637354100199

$ws.Range("A19").Value = "Elopp V1, 600 iterations, Reg=0.2256, RandomSeed=42"
$ws.Range("B19").Value = 0.55000381802544596
$ws.Range("C19").Value = 0.40675764731547498

$ws.Range("A20").Value = "Elopp V1, 50 iterations, Reg=0.2256, RandomSeed=42, HomeAdv=0.27"
$ws.Range("B20").Value = 0.54927045457164503
$ws.Range("C20").Value = 0.406341446076118

$ws.Range("A21").Value = "Elopp V1, 600 iterations, Reg=0.2256, RandomSeed=42, HomeAdv=0.27"
$ws.Range("B21").Value = 0.55033822326860404
$ws.Range("C21").Value = 0.40681785132455101

$ws.Range("A22").Value = "Elopp V1, 100 iterations, Reg=0.2256, RandomSeed=42, HomeAdv=-0.0105"
$ws.Range("B22").Value = 0.54932030697169598
$ws.Range("C22").Value = 0.406441315044622

$ws.Range("A23").Value = "Elopp V1, 600 iterations, Reg=0.2256, RandomSeed=42, HomeAdv=-0.0105"
$ws.Range("B23").Value = 0.54984822782609899
$ws.Range("C23").Value = 0.40670593687463802

# --- Update selection to match target (A24) ---
$ws.Range("A24").Select()
